$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 827.2222
$ws.Range("I55").Value = 189
$ws.Range("K55").Value = 189
$ws.Range("M55").Value = 25
$ws.Range("H111").Value = 14408.479
$ws.Range("I111").Value = 10705.728
$ws.Range("J111").Value = 17802.666
$ws.Range("K111").Value = 32117.184
$ws.Range("L111").Value = 53407.99800000001
$ws.Range("M111").Value = -29050.184
$ws.Range("N111").Value = -59541.99800000001
$ws.Range("H134").Value = 88635.91
$ws.Range("J134").Value = 88635.91
$ws.Range("L134").Value = 88635.91
$ws.Range("N134").Value = -98775.91
$ws.Range("H136").Value = 87996
$ws.Range("J136").Value = 87996
$ws.Range("L136").Value = 87996
$ws.Range("N136").Value = -98196
$ws.Range("H138").Value = 3505.7693
$ws.Range("I138").Value = 1748
$ws.Range("J138").Value = 3706.6572
$ws.Range("K138").Value = 5244
$ws.Range("L138").Value = 11119.9716
$ws.Range("M138").Value = -104
$ws.Range("N138").Value = -21399.9716
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280
$ws.Range("H140").Value = 114492.5
$ws.Range("J140").Value = 114492.5
$ws.Range("L140").Value = 114492.5
$ws.Range("N140").Value = -124852.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3588.6667
$ws.Range("I61").Value = 3979.9375
$ws.Range("K61").Value = 3979.9375
$ws.Range("M61").Value = -3767.9375
$ws.Range("H74").Value = 3335
$ws.Range("I74").Value = 3335
$ws.Range("K74").Value = 3335
$ws.Range("M74").Value = -2461
$ws.Range("H77").Value = 3335
$ws.Range("I77").Value = 3335
$ws.Range("K77").Value = 16675
$ws.Range("M77").Value = -12307
$ws.Range("H133").Value = 79051.8
$ws.Range("J133").Value = 79051.8
$ws.Range("L133").Value = 79051.8
$ws.Range("N133").Value = -84111.8
$ws.Range("H136").Value = 3588.6667
$ws.Range("I136").Value = 3979.9375
$ws.Range("K136").Value = 11939.8125
$ws.Range("M136").Value = -9389.8125
$ws.Range("H137").Value = 79999.5
$ws.Range("J137").Value = 79999
$ws.Range("L137").Value = 79999
$ws.Range("N137").Value = -90199
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 486
$ws.Range("I37").Value = 500
$ws.Range("J37").Value = 481.33334
$ws.Range("K37").Value = 500
$ws.Range("L37").Value = 481.33334
$ws.Range("M37").Value = -363
$ws.Range("N37").Value = -755.33334
$ws.Range("H57").Value = 80425.75
$ws.Range("J57").Value = 79998
$ws.Range("L57").Value = 79998
$ws.Range("N57").Value = -81438
$ws.Range("H136").Value = 80425.75
$ws.Range("J136").Value = 79998
$ws.Range("L136").Value = 79998
$ws.Range("N136").Value = -90198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4389.467
$ws.Range("I58").Value = 4487.28
$ws.Range("K58").Value = 4487.28
$ws.Range("M58").Value = -4284.28
$ws.Range("H136").Value = 4389.467
$ws.Range("I136").Value = 4487.28
$ws.Range("K136").Value = 13461.84
$ws.Range("M136").Value = -10911.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3293.111
$ws.Range("I75").Value = 3049.6667
$ws.Range("K75").Value = 9149.000100000001
$ws.Range("M75").Value = -8151.000100000001
$ws.Range("H78").Value = 3293.111
$ws.Range("I78").Value = 3049.6667
$ws.Range("K78").Value = 27447.0003
$ws.Range("M78").Value = -22455.0003
$ws.Range("H131").Value = 1868.48
$ws.Range("I131").Value = 1320.9333
$ws.Range("K131").Value = 3962.7999
$ws.Range("M131").Value = 1077.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6913.4
$ws.Range("I70").Value = 6765.857
$ws.Range("J70").Value = 7257.6665
$ws.Range("K70").Value = 6765.857
$ws.Range("L70").Value = 7257.6665
$ws.Range("M70").Value = -6495.857
$ws.Range("N70").Value = -7797.6665
$ws.Range("H73").Value = 6913.4
$ws.Range("I73").Value = 6765.857
$ws.Range("J73").Value = 7257.6665
$ws.Range("K73").Value = 6765.857
$ws.Range("L73").Value = 7257.6665
$ws.Range("M73").Value = -5829.857
$ws.Range("N73").Value = -9129.666499999999
$ws.Range("H130").Value = 46600
$ws.Range("I130").Value = 10000
$ws.Range("K130").Value = 10000
$ws.Range("M130").Value = -4980
$ws.Range("H135").Value = 84754.14
$ws.Range("J135").Value = 84754.14
$ws.Range("L135").Value = 84754.14
$ws.Range("N135").Value = -94894.14

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4166.5
$ws.Range("I22").Value = 3000.5
$ws.Range("J22").Value = 4749.5
$ws.Range("K22").Value = 3000.5
$ws.Range("L22").Value = 4749.5
$ws.Range("M22").Value = -2705.5
$ws.Range("N22").Value = -5339.5
$ws.Range("H27").Value = 4166.5
$ws.Range("I27").Value = 3000.5
$ws.Range("J27").Value = 4749.5
$ws.Range("K27").Value = 3000.5
$ws.Range("L27").Value = 4749.5
$ws.Range("M27").Value = -2893.5
$ws.Range("N27").Value = -4963.5
$ws.Range("H122").Value = 9001.799999999999
$ws.Range("I122").Value = 9001.799999999999
$ws.Range("K122").Value = 27005.4
$ws.Range("M122").Value = -24555.4
$ws.Range("H140").Value = 102982.71
$ws.Range("J140").Value = 102982.71
$ws.Range("L140").Value = 102982.71
$ws.Range("N140").Value = -113342.71
$ws.Range("H141").Value = 88652.22
$ws.Range("J141").Value = 88652.22
$ws.Range("L141").Value = 88652.22
$ws.Range("N141").Value = -99012.22

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 4670
$ws.Range("I29").Value = 1010
$ws.Range("K29").Value = 1010
$ws.Range("M29").Value = -720
$ws.Range("H132").Value = 3725.5217
$ws.Range("I132").Value = 3823.0476
$ws.Range("K132").Value = 11469.1428
$ws.Range("M132").Value = -8939.1428
$ws.Range("H136").Value = 3273.0425
$ws.Range("I136").Value = 3175.1316
$ws.Range("K136").Value = 9525.3948
$ws.Range("M136").Value = -6975.3948
$ws.Range("H137").Value = 82373
$ws.Range("J137").Value = 82373
$ws.Range("L137").Value = 82373
$ws.Range("N137").Value = -92573
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
$ws.Range("H141").Value = 68500
$ws.Range("J141").Value = 68500
$ws.Range("L141").Value = 68500
$ws.Range("N141").Value = -78860
